$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.814.70"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.838.40"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.88"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.11"
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  +6.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.23"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0852"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.06"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("D15").Value = "3.281.50"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.977"
$ws.Range("E16").Value = "  +6.11%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.838.42"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "51.876.37"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.48"
$ws.Range("E19").Value = "  +12.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.59"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.88"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.77"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.32"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.62"
$ws.Range("E29").Value = "  +4.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.90"
$ws.Range("E30").Value = "  +7.31%  "
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.32"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.82"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0454"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  +9.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.66"
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.07"
$ws.Range("E38").Value = "  +4.92%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.25"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.20"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").Value = "2.192.04"
$ws.Range("E46").Value = "  +4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.51"
$ws.Range("E47").Value = "  +8.22%  "
$ws.Range("E48").Value = "  +8.08%  "
$ws.Range("E49").Value = "  +22.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.958"
$ws.Range("E50").Value = "  +6.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.55"
$ws.Range("E51").Value = "  +2.58%  "
